# Update New Orleans shard workbook:
#  1. Insert a "State" column into hotel_info between Hotel_Name and City,
#     populated with "Louisiana" for the existing hotel row.
#  2. Reorder the worksheet tabs so review_info comes before hotel_info.

$wb = $excel.ActiveWorkbook

$hotelWs = $wb.Worksheets.Item("hotel_info")
$reviewWs = $wb.Worksheets.Item("review_info")

# Insert a new column C (State) before the existing City column.
$hotelWs.Columns.Item(3).Insert()
$hotelWs.Range("C1").Value = "State"
$hotelWs.Range("C2").Value = "Louisiana"

# Move review_info so it becomes the first sheet (before hotel_info).
$reviewWs.Move($hotelWs)
